$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range to find the last row of data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns D (codeforiati:group-name) and E (codeforiati:group-code) have their
# contents swapped (including the header row), so that D becomes group-code
# and E becomes group-name.
$rangeD = $ws.Range("D1:D$lastRow")
$rangeE = $ws.Range("E1:E$lastRow")

$valuesD = $rangeD.Value()
$valuesE = $rangeE.Value()

$rangeD.Value = $valuesE
$rangeE.Value = $valuesD
